# Apply the "matching" sheet rewrite described by the commit diff.
#
# Sheet "3_" (xl/worksheets/sheet4.xml) is rebuilt: the old prompt row + the
# 6-row term/indicator table (which had a stray orphan row 7) is replaced by
# a tidy 6-row table (1 prompt row + 5 data rows) with the Term/Indicator
# columns swapped (A<->C) and reordered. Sheet "4_" (xl/worksheets/sheet5.xml)
# keeps all of its content - only the remembered cell selection moves.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "3_" : rebuild the matching table
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("3_")

# Drop the old orphan row 7 (just a stray "B" answer with nothing else on
# the row) - this also tightens the sheet dimension down to C6.
$ws3.Rows.Item(7).Delete()

# New prompt text for row 1 (taller, 5-line wrapped prompt).
$ws3.Range("A1").Value = "Here are 5 terms or conditions that appear when a function is defined or called.   Match the function of each term with the term itself."
$ws3.Rows.Item(1).RowHeight = 90

# Row 2: "return"
$ws3.Range("A2").Value = '"return"'
$ws3.Range("B2").Value = "D"
$ws3.Range("C2").Value = 'Indicates the "input" of the function'

# Row 3: "def"
$ws3.Range("A3").Value = '"def"'
$ws3.Range("B3").Value = "E"
$ws3.Range("C3").Value = 'Indicates that the output of the function to the right is being "assigned" to the variable on the left'

# Row 4: Parentheses
$ws3.Range("A4").Value = "Parentheses"
$ws3.Range("B4").Value = "A"
$ws3.Range("C4").Value = "Indicates the lines that define the what the function does to the input"

# Row 5: "=" (equal sign)
$ws3.Range("A5").Value = '"=" (equal sign)'
$ws3.Range("B5").Value = "B"
$ws3.Range("C5").Value = 'Indicates the "output" of the function (and the end of the function definition)'

# Row 6: Indented lines in function definition
$ws3.Range("A6").Value = "Indented lines in function definition"
$ws3.Range("B6").Value = "C"
$ws3.Range("C6").Value = "Indicates the beginning of a function definition"

# ---------------------------------------------------------------------
# Sheet "4_" : content is unchanged, only the saved selection moves
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("4_")

# Restore the saved cell selections, reactivating "3_" then "4_" last so
# "4_" ends up the active tab again (matching the original workbook).
$ws3.Activate() | Out-Null
$ws3.Range("A2").Select() | Out-Null

$ws4.Activate() | Out-Null
$ws4.Range("C13").Select() | Out-Null
